$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the D-column (Price) cells to Text format so numeric-looking strings
# (e.g. "332.34") are stored as text, matching the original inline-string cells,
# instead of being auto-converted to numbers by Excels input parser.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.693.85"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.63"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.34"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +4.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3953"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.08"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08041"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.028"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.06"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.882.64"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.964"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.145"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.13"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001047"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06655"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.18"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.709.10"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.527"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.099.66"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.55"
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.24"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.101"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.593"
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.32"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9740"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09560"
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.339"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06107"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02259"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.236"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.241"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6030"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1903"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.27"
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5697"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.949"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.387"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.50"
$ws.Range("E49").Value = "  +6.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06877"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("E51").Value = "  +14.83%  "
